# Slide 4 ("사용환경 구조" phase diagram): remove the Phase 1-4 rounded-rectangle
# callouts and their "Phase N" text boxes, keeping only the background picture.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -ne "그림 4") {
        $sh.Delete()
    }
}
